$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.602.28'
$ws.Range('E2').Value = '  +2.04%  '
$ws.Range('D3').Value = '2.877.02'
$ws.Range('E3').Value = '  +7.55%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '197.39'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '599.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.554'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.67%  '
$ws.Range('E9').Value = '  +0.59%  '
$ws.Range('D10').Value = '2.871.79'
$ws.Range('E10').Value = '  +7.44%  '
$ws.Range('E11').Value = '  +10.02%  '
$ws.Range('E12').Value = '  -1.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.48%  '
$ws.Range('D14').Value = '3.400.31'
$ws.Range('E14').Value = '  +7.46%  '
$ws.Range('D15').Value = '76.434.68'
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('E16').Value = '  +1.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').Value = '2.863.46'
$ws.Range('E18').Value = '  +6.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('E20').Value = '  +5.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '384.61'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.54%  '
$ws.Range('E22').Value = '  +4.04%  '
$ws.Range('E23').Value = '  +1.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.15%  '
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +2.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.77'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000106'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +11.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.42'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '514.36'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.21%  '
$ws.Range('E34').Value = '  +4.16%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '166.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('E37').Value = '  +4.60%  '
$ws.Range('E38').Value = '  +0.38%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '185.62'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.43%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.346'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.77%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0928'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.45%  '
$ws.Range('E46').Value = '  +4.16%  '
$ws.Range('E47').Value = '  +3.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.580'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.675'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.19%  '
